$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 17 (shifts old rows 17-22 down to 18-23),
#    then copy row 16's formatting/content into it so it matches the
#    original data row's style (borders, fonts, number formats).
$ws.Rows(17).Insert()
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))

# 2. Row 17 becomes the record for CAMILA SAUMETH PALOMINO, keeping her
#    original values (doc/name/period/valor mora/salario basico).
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143466687"
$ws.Range("D17").Value = "CAMILA SAUMETH PALOMINO"
$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 58666
$ws.Range("G17").Value = 2200000
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = ""
$ws.Range("J17").Value = ""

# 3. Row 16 now holds the new worker RAFAEL DE ZUBIRIA CABRALES with new
#    values.
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047468017"
$ws.Range("D16").Value = "RAFAEL DE ZUBIRIA CABRALES"
$ws.Range("E16").Value = "2402"
$ws.Range("F16").Value = 52000
$ws.Range("G16").Value = 1300000

# 4. Update the totals: VALOR MORA is the sum of both workers' mora.
$ws.Range("E11").Value = 110666

# 5. Update worker/period counters now that there are two rows of data.
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2
